$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(94, '$(\langle pq \vert qp \rangle)_{2}$', [double]"0.0004454306859273234"),
    @(20, 'h$_{rs}^{2}$', [double]"0.0003278280035809001"),
    @(63, '$(F_{p}^{\text{SCF}})_{3}$', [double]"0.0002694186150859653"),
    @(81, '$(\langle pq \vert rs \rangle)_{1}$', [double]"0.0002627428358135746"),
    @(46, '$\eta_{s}$', [double]"0.0002021823150789667"),
    @(6, 'h$_{pq}^{2}$', [double]"0.0001878675271866086"),
    @(103, '$(\langle rs\vert rs \rangle)_{3}$', [double]"0.0001708475128863214"),
    @(17, 'h$_{r}^{3}$', [double]"0.0001380112474513354"),
    @(7, 'h$_{pq}^{3}$', [double]"0.000122545963486671"),
    @(102, '$(\langle pq \vert qp \rangle)_{3}$', [double]"0.0001119638415907612"),
    @(21, 'h$_{rs}^{3}$', [double]"0.0001001206740658534"),
    @(12, 'h$_{q}$', [double]"9.858598941052799e-05"),
    @(78, '$(\langle pq \vert qp \rangle)_{0}$', [double]"9.741851548230588e-05"),
    @(74, '$\langle qq \vert qq \rangle$', [double]"9.019656705631658e-05"),
    @(55, '$(F_{p}^{\text{SCF}})_{2}$', [double]"8.868203829582184e-05"),
    @(0, 'h$_{p}^{0}$', [double]"7.670639653803255e-05"),
    @(13, 'h$_{qs}$', [double]"6.985774781381057e-05"),
    @(93, '$(\langle pq \vert pq \rangle)_{2}$', [double]"5.961219675750609e-05"),
    @(89, '$(\langle pq \vert rs \rangle)_{2}$', [double]"5.71895750041908e-05"),
    @(80, '$(\langle rs \vert sr \rangle)_{0}$', [double]"5.511389052411537e-05"),
    @(25, 'typ_2', [double]"5.089233258067151e-05"),
    @(96, '$(\langle rs \vert sr \rangle)_{2}$', [double]"4.966032990948341e-05"),
    @(97, '$(\langle pq \vert rs \rangle)_{3}$', [double]"4.943446548622509e-05"),
    @(101, '$(\langle pq \vert pq \rangle)_{3}$', [double]"4.758798413010874e-05"),
    @(86, '$(\langle pq \vert qp \rangle)_{1}$', [double]"4.111351098361668e-05"),
    @(22, 'h$_{s}$', [double]"3.741824344041914e-05"),
    @(3, 'h$_{p}^{3}$', [double]"3.736542876937551e-05"),
    @(1, 'h$_{p}^{1}$', [double]"3.729346338880552e-05"),
    @(4, 'h$_{pq}^{0}$', [double]"3.709829571482312e-05"),
    @(88, '$(\langle rs \vert sr \rangle)_{1}$', [double]"3.19608406547892e-05"),
    @(2, 'h$_{p}^{2}$', [double]"3.030165225437008e-05"),
    @(50, '$(\eta_{p})_{1}$', [double]"3.008159623205901e-05"),
    @(5, 'h$_{pq}^{1}$', [double]"2.982067715731069e-05"),
    @(104, '$(\langle rs \vert sr \rangle)_{3}$', [double]"2.883206988175777e-05"),
    @(87, '$(\langle rs\vert rs \rangle)_{1}$', [double]"2.774534309261211e-05"),
    @(77, '$(\langle pq \vert pq \rangle)_{0}$', [double]"2.774244701790728e-05"),
    @(75, '$(\langle rr \vert rr \rangle)_{0}$', [double]"2.746035389266966e-05"),
    @(85, '$(\langle pq \vert pq \rangle)_{1}$', [double]"2.449402101285422e-05"),
    @(79, '$(\langle rs\vert rs \rangle)_{0}$', [double]"2.252798411770889e-05"),
    @(10, 'h$_{pr}^{2}$', [double]"2.18140588485958e-05"),
    @(59, '$(F_{r}^{\text{SCF}})_{2}$', [double]"2.177330254275782e-05"),
    @(35, '$F_{q}^{\text{SCF}}$', [double]"2.109251112762774e-05"),
    @(37, '$F_{q}$', [double]"2.027905022754705e-05"),
    @(67, '$(F_{r}^{\text{SCF}})_{3}$', [double]"1.827286740657375e-05"),
    @(73, '$(\langle pp \vert pp \rangle)_{0}$', [double]"1.57854642160131e-05"),
    @(76, '$\langle ss \vert ss \rangle$', [double]"1.570980173486746e-05"),
    @(69, '$(F_{r})_{3}$', [double]"1.545626981063817e-05"),
    @(16, 'h$_{r}^{2}$', [double]"1.545007017337834e-05"),
    @(47, '$(F_{p}^{\text{SCF}})_{1}$', [double]"1.527386852111561e-05"),
    @(11, 'h$_{pr}^{3}$', [double]"1.378914403773468e-05"),
    @(30, 'FA$_{qs}$', [double]"1.330913350417965e-05"),
    @(43, '$F_{s}^{\text{SCF}}$', [double]"1.293126048599399e-05"),
    @(45, '$F_{s}$', [double]"1.234839627410014e-05"),
    @(95, '$(\langle rs\vert rs \rangle)_{2}$', [double]"1.145123799641252e-05"),
    @(8, 'h$_{pr}^{0}$', [double]"1.107144396365168e-05"),
    @(15, 'h$_{r}^{1}$', [double]"1.027571380643737e-05"),
    @(14, 'h$_{r}^{0}$', [double]"1.012816985957754e-05"),
    @(39, '$(F_{r}^{\text{SCF}})_{0}$', [double]"8.526773511377997e-06"),
    @(58, '$(\eta_{p})_{2}$', [double]"8.123893708775933e-06"),
    @(9, 'h$_{pr}^{1}$', [double]"7.674630926620739e-06"),
    @(31, '$(F_{p}^{\text{SCF}})_{0}$', [double]"7.479634687818571e-06"),
    @(70, '$(\eta_{r})_{3}$', [double]"7.179858344393793e-06"),
    @(24, 'typ_1', [double]"7.031517352868992e-06"),
    @(100, '$(\langle rr \vert rr \rangle)_{3}$', [double]"6.922060910896946e-06"),
    @(53, '$(F_{r})_{1}$', [double]"6.26007618308596e-06"),
    @(91, '$(\langle pp \vert pp \rangle)_{2}$', [double]"4.067536048994965e-06"),
    @(71, '$(\langle pq \vert rs \rangle)_{0}$', [double]"3.991355445660034e-06"),
    @(83, '$(\langle pp \vert pp \rangle)_{1}$', [double]"3.79234650949958e-06"),
    @(19, 'h$_{rs}^{1}$', [double]"3.617189580809565e-06"),
    @(18, 'h$_{rs}^{0}$', [double]"3.597900890082835e-06"),
    @(26, 'typ_3', [double]"2.907208461309046e-06"),
    @(92, '$(\langle rr \vert rr \rangle)_{2}$', [double]"2.383408081432526e-06"),
    @(84, '$(\langle rr \vert rr \rangle)_{1}$', [double]"2.378892641810309e-06"),
    @(34, '$(\eta_{p})_{0}$', [double]"1.73833389584705e-06"),
    @(62, '$(\eta_{r})_{2}$', [double]"1.603961242501242e-06"),
    @(41, '$(F_{r})_{0}$', [double]"9.912576064833218e-07"),
    @(51, '$(F_{r}^{\text{SCF}})_{1}$', [double]"9.231544017611491e-07"),
    @(57, '$(F_{p})_{2}$', [double]"9.226245036920213e-07"),
    @(98, '$(\langle pq \vert sr \rangle)_{3}$', [double]"2.47336876731498e-07"),
    @(29, 'FI$_{qs}$', [double]"1.826275964185884e-07"),
    @(49, '$(F_{p})_{1}$', [double]"6.199446060321493e-08"),
    @(72, '$(\langle pq \vert sr \rangle)_{0}$', [double]"5.767702474039683e-08"),
    @(90, '$(\langle pq \vert sr \rangle)_{2}$', [double]"5.687813813972815e-08"),
    @(82, '$(\langle pq \vert sr \rangle)_{1}$', [double]"5.565174244116233e-08"),
    @(99, '$(\langle pp \vert pp \rangle)_{3}$', [double]"5.447333865161393e-08"),
    @(33, '$(F_{p})_{0}$', [double]"5.17947966401994e-08"),
    @(61, '$(F_{r})_{2}$', [double]"5.139361212872752e-08"),
    @(65, '$(F_{p})_{3}$', [double]"5.130020151627852e-08"),
    @(54, '$(\eta_{r})_{1}$', [double]"4.935108610399632e-08"),
    @(42, '$(\eta_{r})_{0}$', [double]"4.91458114053784e-08"),
    @(48, '$(\omega_{p})_{1}$', [double]"3.912303268612381e-08"),
    @(38, '$\eta_{q}$', [double]"3.811554435581008e-08"),
    @(52, '$(\omega_{r})_{1}$', [double]"3.739483368226801e-08"),
    @(40, '$(\omega_{r})_{0}$', [double]"3.338282112362755e-08"),
    @(56, '$(\omega_{p})_{2}$', [double]"2.923466032638399e-08"),
    @(66, '$(\eta_{p})_{3}$', [double]"2.653840608747067e-08"),
    @(23, 'typ_0', [double]"2.65182405928935e-08"),
    @(32, '$(\omega_{p})_{0}$', [double]"1.669821466401023e-08"),
    @(60, '$(\omega_{r})_{2}$', [double]"8.419080239005867e-09"),
    @(27, '$\mathbf{b}$', [double]"8.133954095740331e-09"),
    @(64, '$(\omega_{p})_{3}$', [double]"8.063172763450642e-09"),
    @(68, '$(\omega_{r})_{3}$', [double]"7.368139515703121e-09"),
    @(44, '$\omega_{s}$', [double]"6.505880660877509e-09"),
    @(28, 'F$_{qs}$', [double]"1.306780828268086e-09"),
    @(36, '$\omega_{q}$', [double]"5.891473602562746e-10")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
